$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Snapshot original row data (values/formulas) before moving things around
$snapshot = @{}
$snapshot[10] = @{}
$snapshot[10][1] = @{ Kind = "String"; Val = "A 60406-2025" }
$snapshot[10][2] = @{ Kind = "Number"; Val = 45995.43050925926 }
$snapshot[10][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[10][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[10][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[10][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[10][7] = @{ Kind = "Number"; Val = 2.3 }
$snapshot[10][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][9] = @{ Kind = "Number"; Val = 2 }
$snapshot[10][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[10][17] = @{ Kind = "Number"; Val = 2 }
$snapshot[10][18] = @{ Kind = "String"; Val = "Björksplintborre`r`nMindre märgborre" }
$snapshot[10][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 60406-2025 artfynd.xlsx`", `"A 60406-2025`")" }
$snapshot[10][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 60406-2025 karta.png`", `"A 60406-2025`")" }
$snapshot[10][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 60406-2025 FSC-klagomål.docx`", `"A 60406-2025`")" }
$snapshot[10][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 60406-2025 FSC-klagomål mail.docx`", `"A 60406-2025`")" }
$snapshot[10][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 60406-2025 tillsynsbegäran.docx`", `"A 60406-2025`")" }
$snapshot[10][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 60406-2025 tillsynsbegäran mail.docx`", `"A 60406-2025`")" }
$snapshot[11] = @{}
$snapshot[11][1] = @{ Kind = "String"; Val = "A 1697-2023" }
$snapshot[11][2] = @{ Kind = "Number"; Val = 44938.0 }
$snapshot[11][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[11][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[11][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[11][7] = @{ Kind = "Number"; Val = 1.1 }
$snapshot[11][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][9] = @{ Kind = "Number"; Val = 1 }
$snapshot[11][10] = @{ Kind = "Number"; Val = 1 }
$snapshot[11][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][15] = @{ Kind = "Number"; Val = 1 }
$snapshot[11][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[11][17] = @{ Kind = "Number"; Val = 2 }
$snapshot[11][18] = @{ Kind = "String"; Val = "Lunglav`r`nBårdlav" }
$snapshot[11][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 1697-2023 artfynd.xlsx`", `"A 1697-2023`")" }
$snapshot[11][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 1697-2023 karta.png`", `"A 1697-2023`")" }
$snapshot[11][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 1697-2023 FSC-klagomål.docx`", `"A 1697-2023`")" }
$snapshot[11][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 1697-2023 FSC-klagomål mail.docx`", `"A 1697-2023`")" }
$snapshot[11][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 1697-2023 tillsynsbegäran.docx`", `"A 1697-2023`")" }
$snapshot[11][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 1697-2023 tillsynsbegäran mail.docx`", `"A 1697-2023`")" }
$snapshot[13] = @{}
$snapshot[13][1] = @{ Kind = "String"; Val = "A 30067-2025" }
$snapshot[13][2] = @{ Kind = "Number"; Val = 45826.0 }
$snapshot[13][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[13][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[13][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[13][7] = @{ Kind = "Number"; Val = 11.9 }
$snapshot[13][8] = @{ Kind = "Number"; Val = 2 }
$snapshot[13][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[13][17] = @{ Kind = "Number"; Val = 2 }
$snapshot[13][18] = @{ Kind = "String"; Val = "Åkergroda`r`nVanlig groda" }
$snapshot[13][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 30067-2025 artfynd.xlsx`", `"A 30067-2025`")" }
$snapshot[13][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 30067-2025 karta.png`", `"A 30067-2025`")" }
$snapshot[13][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 30067-2025 FSC-klagomål.docx`", `"A 30067-2025`")" }
$snapshot[13][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 30067-2025 FSC-klagomål mail.docx`", `"A 30067-2025`")" }
$snapshot[13][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 30067-2025 tillsynsbegäran.docx`", `"A 30067-2025`")" }
$snapshot[13][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 30067-2025 tillsynsbegäran mail.docx`", `"A 30067-2025`")" }
$snapshot[17] = @{}
$snapshot[17][1] = @{ Kind = "String"; Val = "A 1691-2023" }
$snapshot[17][2] = @{ Kind = "Number"; Val = 44938.0 }
$snapshot[17][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[17][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[17][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[17][7] = @{ Kind = "Number"; Val = 1.7 }
$snapshot[17][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][9] = @{ Kind = "Number"; Val = 1 }
$snapshot[17][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[17][17] = @{ Kind = "Number"; Val = 1 }
$snapshot[17][18] = @{ Kind = "String"; Val = "Tibast" }
$snapshot[17][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 1691-2023 artfynd.xlsx`", `"A 1691-2023`")" }
$snapshot[17][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 1691-2023 karta.png`", `"A 1691-2023`")" }
$snapshot[17][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 1691-2023 FSC-klagomål.docx`", `"A 1691-2023`")" }
$snapshot[17][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 1691-2023 FSC-klagomål mail.docx`", `"A 1691-2023`")" }
$snapshot[17][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 1691-2023 tillsynsbegäran.docx`", `"A 1691-2023`")" }
$snapshot[17][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 1691-2023 tillsynsbegäran mail.docx`", `"A 1691-2023`")" }
$snapshot[19] = @{}
$snapshot[19][1] = @{ Kind = "String"; Val = "A 18100-2022" }
$snapshot[19][2] = @{ Kind = "Number"; Val = 44684.0 }
$snapshot[19][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[19][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[19][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[19][7] = @{ Kind = "Number"; Val = 3.1 }
$snapshot[19][8] = @{ Kind = "Number"; Val = 1 }
$snapshot[19][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[19][17] = @{ Kind = "Number"; Val = 1 }
$snapshot[19][18] = @{ Kind = "String"; Val = "Gullviva" }
$snapshot[19][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 18100-2022 artfynd.xlsx`", `"A 18100-2022`")" }
$snapshot[19][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 18100-2022 karta.png`", `"A 18100-2022`")" }
$snapshot[19][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 18100-2022 FSC-klagomål.docx`", `"A 18100-2022`")" }
$snapshot[19][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 18100-2022 FSC-klagomål mail.docx`", `"A 18100-2022`")" }
$snapshot[19][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 18100-2022 tillsynsbegäran.docx`", `"A 18100-2022`")" }
$snapshot[19][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 18100-2022 tillsynsbegäran mail.docx`", `"A 18100-2022`")" }
$snapshot[20] = @{}
$snapshot[20][1] = @{ Kind = "String"; Val = "A 33801-2025" }
$snapshot[20][2] = @{ Kind = "Number"; Val = 45842.0 }
$snapshot[20][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[20][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[20][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[20][7] = @{ Kind = "Number"; Val = 1.3 }
$snapshot[20][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][9] = @{ Kind = "Number"; Val = 1 }
$snapshot[20][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[20][17] = @{ Kind = "Number"; Val = 1 }
$snapshot[20][18] = @{ Kind = "String"; Val = "Strutbräken" }
$snapshot[20][19] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/artfynd/A 33801-2025 artfynd.xlsx`", `"A 33801-2025`")" }
$snapshot[20][20] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/kartor/A 33801-2025 karta.png`", `"A 33801-2025`")" }
$snapshot[20][22] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomål/A 33801-2025 FSC-klagomål.docx`", `"A 33801-2025`")" }
$snapshot[20][23] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/klagomålsmail/A 33801-2025 FSC-klagomål mail.docx`", `"A 33801-2025`")" }
$snapshot[20][24] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsyn/A 33801-2025 tillsynsbegäran.docx`", `"A 33801-2025`")" }
$snapshot[20][25] = @{ Kind = "Formula"; Val = "=HYPERLINK(`"https://klasma.github.io/Logging_0127/tillsynsmail/A 33801-2025 tillsynsbegäran mail.docx`", `"A 33801-2025`")" }
$snapshot[34] = @{}
$snapshot[34][1] = @{ Kind = "String"; Val = "A 11867-2024" }
$snapshot[34][2] = @{ Kind = "Number"; Val = 45376.45421296296 }
$snapshot[34][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[34][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[34][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[34][7] = @{ Kind = "Number"; Val = 2.9 }
$snapshot[34][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[34][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[35] = @{}
$snapshot[35][1] = @{ Kind = "String"; Val = "A 48656-2022" }
$snapshot[35][2] = @{ Kind = "Number"; Val = 44859.439930555556 }
$snapshot[35][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[35][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[35][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[35][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[35][7] = @{ Kind = "Number"; Val = 1.1 }
$snapshot[35][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[35][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[36] = @{}
$snapshot[36][1] = @{ Kind = "String"; Val = "A 16149-2024" }
$snapshot[36][2] = @{ Kind = "Number"; Val = 45406.56984953704 }
$snapshot[36][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[36][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[36][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[36][7] = @{ Kind = "Number"; Val = 1 }
$snapshot[36][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[36][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[37] = @{}
$snapshot[37][1] = @{ Kind = "String"; Val = "A 16145-2024" }
$snapshot[37][2] = @{ Kind = "Number"; Val = 45406.565613425926 }
$snapshot[37][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[37][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[37][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[37][7] = @{ Kind = "Number"; Val = 0.6 }
$snapshot[37][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[37][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[38] = @{}
$snapshot[38][1] = @{ Kind = "String"; Val = "A 46405-2022" }
$snapshot[38][2] = @{ Kind = "Number"; Val = 44846.0 }
$snapshot[38][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[38][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[38][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[38][7] = @{ Kind = "Number"; Val = 3.4 }
$snapshot[38][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[38][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[39] = @{}
$snapshot[39][1] = @{ Kind = "String"; Val = "A 37095-2021" }
$snapshot[39][2] = @{ Kind = "Number"; Val = 44395.0 }
$snapshot[39][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[39][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[39][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[39][7] = @{ Kind = "Number"; Val = 3.5 }
$snapshot[39][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[39][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[40] = @{}
$snapshot[40][1] = @{ Kind = "String"; Val = "A 38194-2022" }
$snapshot[40][2] = @{ Kind = "Number"; Val = 44812.0 }
$snapshot[40][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[40][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[40][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[40][7] = @{ Kind = "Number"; Val = 1.6 }
$snapshot[40][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[40][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[41] = @{}
$snapshot[41][1] = @{ Kind = "String"; Val = "A 28983-2024" }
$snapshot[41][2] = @{ Kind = "Number"; Val = 45481.477314814816 }
$snapshot[41][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[41][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[41][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[41][7] = @{ Kind = "Number"; Val = 0.8 }
$snapshot[41][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[41][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[42] = @{}
$snapshot[42][1] = @{ Kind = "String"; Val = "A 45945-2025" }
$snapshot[42][2] = @{ Kind = "Number"; Val = 45924.0 }
$snapshot[42][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[42][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[42][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[42][7] = @{ Kind = "Number"; Val = 3.1 }
$snapshot[42][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[42][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[43] = @{}
$snapshot[43][1] = @{ Kind = "String"; Val = "A 73151-2021" }
$snapshot[43][2] = @{ Kind = "Number"; Val = 44550.5162037037 }
$snapshot[43][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[43][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[43][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[43][7] = @{ Kind = "Number"; Val = 0.5 }
$snapshot[43][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[43][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[44] = @{}
$snapshot[44][1] = @{ Kind = "String"; Val = "A 73194-2021" }
$snapshot[44][2] = @{ Kind = "Number"; Val = 44550.61508101852 }
$snapshot[44][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[44][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[44][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[44][7] = @{ Kind = "Number"; Val = 1.1 }
$snapshot[44][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[44][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[45] = @{}
$snapshot[45][1] = @{ Kind = "String"; Val = "A 73155-2021" }
$snapshot[45][2] = @{ Kind = "Number"; Val = 44550.0 }
$snapshot[45][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[45][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[45][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[45][7] = @{ Kind = "Number"; Val = 0.8 }
$snapshot[45][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[45][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[46] = @{}
$snapshot[46][1] = @{ Kind = "String"; Val = "A 30378-2024" }
$snapshot[46][2] = @{ Kind = "Number"; Val = 45491.0 }
$snapshot[46][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[46][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[46][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[46][7] = @{ Kind = "Number"; Val = 0.3 }
$snapshot[46][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[46][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[48] = @{}
$snapshot[48][1] = @{ Kind = "String"; Val = "A 15905-2022" }
$snapshot[48][2] = @{ Kind = "Number"; Val = 44664.0 }
$snapshot[48][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[48][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[48][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[48][7] = @{ Kind = "Number"; Val = 11.2 }
$snapshot[48][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[48][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[50] = @{}
$snapshot[50][1] = @{ Kind = "String"; Val = "A 61514-2022" }
$snapshot[50][2] = @{ Kind = "Number"; Val = 44916.0 }
$snapshot[50][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[50][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[50][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[50][7] = @{ Kind = "Number"; Val = 1.3 }
$snapshot[50][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[50][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[51] = @{}
$snapshot[51][1] = @{ Kind = "String"; Val = "A 15692-2023" }
$snapshot[51][2] = @{ Kind = "Number"; Val = 45020.0 }
$snapshot[51][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[51][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[51][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[51][7] = @{ Kind = "Number"; Val = 2.1 }
$snapshot[51][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[51][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[52] = @{}
$snapshot[52][1] = @{ Kind = "String"; Val = "A 10815-2021" }
$snapshot[52][2] = @{ Kind = "Number"; Val = 44259.0 }
$snapshot[52][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[52][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[52][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[52][7] = @{ Kind = "Number"; Val = 1.6 }
$snapshot[52][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[52][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[53] = @{}
$snapshot[53][1] = @{ Kind = "String"; Val = "A 14202-2025" }
$snapshot[53][2] = @{ Kind = "Number"; Val = 45740.57461805556 }
$snapshot[53][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[53][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[53][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[53][7] = @{ Kind = "Number"; Val = 0.9 }
$snapshot[53][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[53][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[55] = @{}
$snapshot[55][1] = @{ Kind = "String"; Val = "A 21285-2023" }
$snapshot[55][2] = @{ Kind = "Number"; Val = 45062.0 }
$snapshot[55][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[55][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[55][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[55][7] = @{ Kind = "Number"; Val = 0.8 }
$snapshot[55][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[55][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[56] = @{}
$snapshot[56][1] = @{ Kind = "String"; Val = "A 15690-2023" }
$snapshot[56][2] = @{ Kind = "Number"; Val = 45020.0 }
$snapshot[56][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[56][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[56][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[56][7] = @{ Kind = "Number"; Val = 5 }
$snapshot[56][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[56][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[57] = @{}
$snapshot[57][1] = @{ Kind = "String"; Val = "A 3021-2025" }
$snapshot[57][2] = @{ Kind = "Number"; Val = 45678.57962962963 }
$snapshot[57][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[57][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[57][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[57][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[57][7] = @{ Kind = "Number"; Val = 1.4 }
$snapshot[57][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[57][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[58] = @{}
$snapshot[58][1] = @{ Kind = "String"; Val = "A 36097-2021" }
$snapshot[58][2] = @{ Kind = "Number"; Val = 44389.0 }
$snapshot[58][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[58][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[58][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[58][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[58][7] = @{ Kind = "Number"; Val = 2.5 }
$snapshot[58][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[58][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[59] = @{}
$snapshot[59][1] = @{ Kind = "String"; Val = "A 20649-2022" }
$snapshot[59][2] = @{ Kind = "Number"; Val = 44700.0 }
$snapshot[59][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[59][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[59][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[59][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[59][7] = @{ Kind = "Number"; Val = 6 }
$snapshot[59][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[59][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[60] = @{}
$snapshot[60][1] = @{ Kind = "String"; Val = "A 6059-2022" }
$snapshot[60][2] = @{ Kind = "Number"; Val = 44599.0 }
$snapshot[60][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[60][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[60][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[60][7] = @{ Kind = "Number"; Val = 1.1 }
$snapshot[60][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[60][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[61] = @{}
$snapshot[61][1] = @{ Kind = "String"; Val = "A 61225-2022" }
$snapshot[61][2] = @{ Kind = "Number"; Val = 44915.0 }
$snapshot[61][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[61][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[61][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[61][7] = @{ Kind = "Number"; Val = 1.9 }
$snapshot[61][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[61][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[62] = @{}
$snapshot[62][1] = @{ Kind = "String"; Val = "A 8528-2024" }
$snapshot[62][2] = @{ Kind = "Number"; Val = 45355.49099537037 }
$snapshot[62][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[62][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[62][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[62][7] = @{ Kind = "Number"; Val = 0.4 }
$snapshot[62][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[62][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[63] = @{}
$snapshot[63][1] = @{ Kind = "String"; Val = "A 3686-2022" }
$snapshot[63][2] = @{ Kind = "Number"; Val = 44586.57239583333 }
$snapshot[63][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[63][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[63][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[63][7] = @{ Kind = "Number"; Val = 3.8 }
$snapshot[63][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[63][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[64] = @{}
$snapshot[64][1] = @{ Kind = "String"; Val = "A 35197-2024" }
$snapshot[64][2] = @{ Kind = "Number"; Val = 45530.0 }
$snapshot[64][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[64][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[64][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[64][7] = @{ Kind = "Number"; Val = 2.6 }
$snapshot[64][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[64][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[65] = @{}
$snapshot[65][1] = @{ Kind = "String"; Val = "A 9032-2023" }
$snapshot[65][2] = @{ Kind = "Number"; Val = 44979.0 }
$snapshot[65][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[65][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[65][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[65][7] = @{ Kind = "Number"; Val = 2.6 }
$snapshot[65][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[65][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[66] = @{}
$snapshot[66][1] = @{ Kind = "String"; Val = "A 25217-2024" }
$snapshot[66][2] = @{ Kind = "Number"; Val = 45462.0 }
$snapshot[66][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[66][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[66][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[66][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[66][7] = @{ Kind = "Number"; Val = 2 }
$snapshot[66][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[66][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[67] = @{}
$snapshot[67][1] = @{ Kind = "String"; Val = "A 16483-2023" }
$snapshot[67][2] = @{ Kind = "Number"; Val = 45029.0 }
$snapshot[67][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[67][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[67][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[67][7] = @{ Kind = "Number"; Val = 3.1 }
$snapshot[67][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[67][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[68] = @{}
$snapshot[68][1] = @{ Kind = "String"; Val = "A 30622-2024" }
$snapshot[68][2] = @{ Kind = "Number"; Val = 45495.0 }
$snapshot[68][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[68][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[68][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[68][7] = @{ Kind = "Number"; Val = 0.3 }
$snapshot[68][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[68][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[69] = @{}
$snapshot[69][1] = @{ Kind = "String"; Val = "A 61320-2024" }
$snapshot[69][2] = @{ Kind = "Number"; Val = 45646.33262731481 }
$snapshot[69][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[69][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[69][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[69][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[69][7] = @{ Kind = "Number"; Val = 0.9 }
$snapshot[69][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[69][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[70] = @{}
$snapshot[70][1] = @{ Kind = "String"; Val = "A 30925-2023" }
$snapshot[70][2] = @{ Kind = "Number"; Val = 45113.0 }
$snapshot[70][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[70][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[70][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[70][6] = @{ Kind = "String"; Val = "Övriga Aktiebolag" }
$snapshot[70][7] = @{ Kind = "Number"; Val = 0.8 }
$snapshot[70][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[70][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[71] = @{}
$snapshot[71][1] = @{ Kind = "String"; Val = "A 45869-2022" }
$snapshot[71][2] = @{ Kind = "Number"; Val = 44844.0 }
$snapshot[71][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[71][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[71][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[71][7] = @{ Kind = "Number"; Val = 3.4 }
$snapshot[71][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[71][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[72] = @{}
$snapshot[72][1] = @{ Kind = "String"; Val = "A 14217-2025" }
$snapshot[72][2] = @{ Kind = "Number"; Val = 45740.58534722222 }
$snapshot[72][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[72][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[72][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[72][7] = @{ Kind = "Number"; Val = 3.4 }
$snapshot[72][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[72][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[73] = @{}
$snapshot[73][1] = @{ Kind = "String"; Val = "A 61242-2022" }
$snapshot[73][2] = @{ Kind = "Number"; Val = 44915.0 }
$snapshot[73][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[73][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[73][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[73][7] = @{ Kind = "Number"; Val = 2.8 }
$snapshot[73][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[73][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[74] = @{}
$snapshot[74][1] = @{ Kind = "String"; Val = "A 11876-2024" }
$snapshot[74][2] = @{ Kind = "Number"; Val = 45376.470729166664 }
$snapshot[74][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[74][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[74][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[74][7] = @{ Kind = "Number"; Val = 1 }
$snapshot[74][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[74][17] = @{ Kind = "Number"; Val = 0 }
$snapshot[75] = @{}
$snapshot[75][1] = @{ Kind = "String"; Val = "A 16142-2024" }
$snapshot[75][2] = @{ Kind = "Number"; Val = 45406.0 }
$snapshot[75][3] = @{ Kind = "Number"; Val = 46059.0 }
$snapshot[75][4] = @{ Kind = "String"; Val = "STOCKHOLMS LÄN" }
$snapshot[75][5] = @{ Kind = "String"; Val = "BOTKYRKA" }
$snapshot[75][7] = @{ Kind = "Number"; Val = 1.5 }
$snapshot[75][8] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][9] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][10] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][11] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][12] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][13] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][14] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][15] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][16] = @{ Kind = "Number"; Val = 0 }
$snapshot[75][17] = @{ Kind = "Number"; Val = 0 }

# Apply snapshot data to target rows, and clear columns that should become empty
# Target row 10 <- source row 11
$cell = $ws.Cells.Item(10, 1)
$entry = $snapshot[11][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 2)
$entry = $snapshot[11][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 3)
$entry = $snapshot[11][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 4)
$entry = $snapshot[11][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 5)
$entry = $snapshot[11][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 7)
$entry = $snapshot[11][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 8)
$entry = $snapshot[11][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 9)
$entry = $snapshot[11][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 10)
$entry = $snapshot[11][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 11)
$entry = $snapshot[11][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 12)
$entry = $snapshot[11][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 13)
$entry = $snapshot[11][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 14)
$entry = $snapshot[11][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 15)
$entry = $snapshot[11][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 16)
$entry = $snapshot[11][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 17)
$entry = $snapshot[11][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 18)
$entry = $snapshot[11][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 19)
$entry = $snapshot[11][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 20)
$entry = $snapshot[11][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 22)
$entry = $snapshot[11][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 23)
$entry = $snapshot[11][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 24)
$entry = $snapshot[11][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(10, 25)
$entry = $snapshot[11][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(10, 6).ClearContents()

# Target row 11 <- source row 13
$cell = $ws.Cells.Item(11, 1)
$entry = $snapshot[13][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 2)
$entry = $snapshot[13][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 3)
$entry = $snapshot[13][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 4)
$entry = $snapshot[13][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 5)
$entry = $snapshot[13][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 7)
$entry = $snapshot[13][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 8)
$entry = $snapshot[13][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 9)
$entry = $snapshot[13][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 10)
$entry = $snapshot[13][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 11)
$entry = $snapshot[13][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 12)
$entry = $snapshot[13][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 13)
$entry = $snapshot[13][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 14)
$entry = $snapshot[13][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 15)
$entry = $snapshot[13][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 16)
$entry = $snapshot[13][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 17)
$entry = $snapshot[13][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 18)
$entry = $snapshot[13][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 19)
$entry = $snapshot[13][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 20)
$entry = $snapshot[13][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 22)
$entry = $snapshot[13][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 23)
$entry = $snapshot[13][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 24)
$entry = $snapshot[13][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(11, 25)
$entry = $snapshot[13][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 13 <- source row 10
$cell = $ws.Cells.Item(13, 1)
$entry = $snapshot[10][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 2)
$entry = $snapshot[10][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 3)
$entry = $snapshot[10][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 4)
$entry = $snapshot[10][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 5)
$entry = $snapshot[10][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 6)
$entry = $snapshot[10][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 7)
$entry = $snapshot[10][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 8)
$entry = $snapshot[10][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 9)
$entry = $snapshot[10][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 10)
$entry = $snapshot[10][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 11)
$entry = $snapshot[10][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 12)
$entry = $snapshot[10][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 13)
$entry = $snapshot[10][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 14)
$entry = $snapshot[10][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 15)
$entry = $snapshot[10][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 16)
$entry = $snapshot[10][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 17)
$entry = $snapshot[10][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 18)
$entry = $snapshot[10][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 19)
$entry = $snapshot[10][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 20)
$entry = $snapshot[10][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 22)
$entry = $snapshot[10][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 23)
$entry = $snapshot[10][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 24)
$entry = $snapshot[10][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(13, 25)
$entry = $snapshot[10][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 17 <- source row 20
$cell = $ws.Cells.Item(17, 1)
$entry = $snapshot[20][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 2)
$entry = $snapshot[20][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 3)
$entry = $snapshot[20][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 4)
$entry = $snapshot[20][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 5)
$entry = $snapshot[20][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 7)
$entry = $snapshot[20][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 8)
$entry = $snapshot[20][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 9)
$entry = $snapshot[20][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 10)
$entry = $snapshot[20][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 11)
$entry = $snapshot[20][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 12)
$entry = $snapshot[20][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 13)
$entry = $snapshot[20][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 14)
$entry = $snapshot[20][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 15)
$entry = $snapshot[20][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 16)
$entry = $snapshot[20][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 17)
$entry = $snapshot[20][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 18)
$entry = $snapshot[20][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 19)
$entry = $snapshot[20][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 20)
$entry = $snapshot[20][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 22)
$entry = $snapshot[20][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 23)
$entry = $snapshot[20][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 24)
$entry = $snapshot[20][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(17, 25)
$entry = $snapshot[20][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 19 <- source row 17
$cell = $ws.Cells.Item(19, 1)
$entry = $snapshot[17][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 2)
$entry = $snapshot[17][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 3)
$entry = $snapshot[17][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 4)
$entry = $snapshot[17][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 5)
$entry = $snapshot[17][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 7)
$entry = $snapshot[17][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 8)
$entry = $snapshot[17][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 9)
$entry = $snapshot[17][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 10)
$entry = $snapshot[17][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 11)
$entry = $snapshot[17][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 12)
$entry = $snapshot[17][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 13)
$entry = $snapshot[17][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 14)
$entry = $snapshot[17][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 15)
$entry = $snapshot[17][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 16)
$entry = $snapshot[17][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 17)
$entry = $snapshot[17][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 18)
$entry = $snapshot[17][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 19)
$entry = $snapshot[17][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 20)
$entry = $snapshot[17][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 22)
$entry = $snapshot[17][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 23)
$entry = $snapshot[17][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 24)
$entry = $snapshot[17][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(19, 25)
$entry = $snapshot[17][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 20 <- source row 19
$cell = $ws.Cells.Item(20, 1)
$entry = $snapshot[19][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 2)
$entry = $snapshot[19][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 3)
$entry = $snapshot[19][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 4)
$entry = $snapshot[19][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 5)
$entry = $snapshot[19][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 7)
$entry = $snapshot[19][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 8)
$entry = $snapshot[19][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 9)
$entry = $snapshot[19][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 10)
$entry = $snapshot[19][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 11)
$entry = $snapshot[19][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 12)
$entry = $snapshot[19][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 13)
$entry = $snapshot[19][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 14)
$entry = $snapshot[19][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 15)
$entry = $snapshot[19][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 16)
$entry = $snapshot[19][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 17)
$entry = $snapshot[19][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 18)
$entry = $snapshot[19][18]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 19)
$entry = $snapshot[19][19]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 20)
$entry = $snapshot[19][20]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 22)
$entry = $snapshot[19][22]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 23)
$entry = $snapshot[19][23]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 24)
$entry = $snapshot[19][24]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(20, 25)
$entry = $snapshot[19][25]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 34 <- source row 53
$cell = $ws.Cells.Item(34, 1)
$entry = $snapshot[53][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 2)
$entry = $snapshot[53][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 3)
$entry = $snapshot[53][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 4)
$entry = $snapshot[53][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 5)
$entry = $snapshot[53][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 7)
$entry = $snapshot[53][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 8)
$entry = $snapshot[53][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 9)
$entry = $snapshot[53][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 10)
$entry = $snapshot[53][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 11)
$entry = $snapshot[53][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 12)
$entry = $snapshot[53][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 13)
$entry = $snapshot[53][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 14)
$entry = $snapshot[53][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 15)
$entry = $snapshot[53][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 16)
$entry = $snapshot[53][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(34, 17)
$entry = $snapshot[53][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 35 <- source row 73
$cell = $ws.Cells.Item(35, 1)
$entry = $snapshot[73][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 2)
$entry = $snapshot[73][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 3)
$entry = $snapshot[73][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 4)
$entry = $snapshot[73][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 5)
$entry = $snapshot[73][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 7)
$entry = $snapshot[73][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 8)
$entry = $snapshot[73][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 9)
$entry = $snapshot[73][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 10)
$entry = $snapshot[73][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 11)
$entry = $snapshot[73][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 12)
$entry = $snapshot[73][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 13)
$entry = $snapshot[73][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 14)
$entry = $snapshot[73][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 15)
$entry = $snapshot[73][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 16)
$entry = $snapshot[73][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(35, 17)
$entry = $snapshot[73][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(35, 6).ClearContents()

# Target row 36 <- source row 55
$cell = $ws.Cells.Item(36, 1)
$entry = $snapshot[55][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 2)
$entry = $snapshot[55][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 3)
$entry = $snapshot[55][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 4)
$entry = $snapshot[55][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 5)
$entry = $snapshot[55][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 7)
$entry = $snapshot[55][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 8)
$entry = $snapshot[55][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 9)
$entry = $snapshot[55][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 10)
$entry = $snapshot[55][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 11)
$entry = $snapshot[55][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 12)
$entry = $snapshot[55][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 13)
$entry = $snapshot[55][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 14)
$entry = $snapshot[55][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 15)
$entry = $snapshot[55][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 16)
$entry = $snapshot[55][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(36, 17)
$entry = $snapshot[55][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 37 <- source row 70
$cell = $ws.Cells.Item(37, 1)
$entry = $snapshot[70][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 2)
$entry = $snapshot[70][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 3)
$entry = $snapshot[70][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 4)
$entry = $snapshot[70][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 5)
$entry = $snapshot[70][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 6)
$entry = $snapshot[70][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 7)
$entry = $snapshot[70][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 8)
$entry = $snapshot[70][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 9)
$entry = $snapshot[70][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 10)
$entry = $snapshot[70][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 11)
$entry = $snapshot[70][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 12)
$entry = $snapshot[70][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 13)
$entry = $snapshot[70][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 14)
$entry = $snapshot[70][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 15)
$entry = $snapshot[70][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 16)
$entry = $snapshot[70][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(37, 17)
$entry = $snapshot[70][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 38 <- source row 56
$cell = $ws.Cells.Item(38, 1)
$entry = $snapshot[56][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 2)
$entry = $snapshot[56][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 3)
$entry = $snapshot[56][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 4)
$entry = $snapshot[56][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 5)
$entry = $snapshot[56][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 7)
$entry = $snapshot[56][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 8)
$entry = $snapshot[56][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 9)
$entry = $snapshot[56][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 10)
$entry = $snapshot[56][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 11)
$entry = $snapshot[56][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 12)
$entry = $snapshot[56][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 13)
$entry = $snapshot[56][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 14)
$entry = $snapshot[56][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 15)
$entry = $snapshot[56][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 16)
$entry = $snapshot[56][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(38, 17)
$entry = $snapshot[56][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 39 <- source row 72
$cell = $ws.Cells.Item(39, 1)
$entry = $snapshot[72][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 2)
$entry = $snapshot[72][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 3)
$entry = $snapshot[72][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 4)
$entry = $snapshot[72][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 5)
$entry = $snapshot[72][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 7)
$entry = $snapshot[72][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 8)
$entry = $snapshot[72][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 9)
$entry = $snapshot[72][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 10)
$entry = $snapshot[72][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 11)
$entry = $snapshot[72][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 12)
$entry = $snapshot[72][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 13)
$entry = $snapshot[72][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 14)
$entry = $snapshot[72][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 15)
$entry = $snapshot[72][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 16)
$entry = $snapshot[72][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(39, 17)
$entry = $snapshot[72][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 40 <- source row 46
$cell = $ws.Cells.Item(40, 1)
$entry = $snapshot[46][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 2)
$entry = $snapshot[46][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 3)
$entry = $snapshot[46][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 4)
$entry = $snapshot[46][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 5)
$entry = $snapshot[46][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 7)
$entry = $snapshot[46][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 8)
$entry = $snapshot[46][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 9)
$entry = $snapshot[46][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 10)
$entry = $snapshot[46][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 11)
$entry = $snapshot[46][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 12)
$entry = $snapshot[46][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 13)
$entry = $snapshot[46][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 14)
$entry = $snapshot[46][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 15)
$entry = $snapshot[46][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 16)
$entry = $snapshot[46][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(40, 17)
$entry = $snapshot[46][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 41 <- source row 69
$cell = $ws.Cells.Item(41, 1)
$entry = $snapshot[69][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 2)
$entry = $snapshot[69][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 3)
$entry = $snapshot[69][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 4)
$entry = $snapshot[69][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 5)
$entry = $snapshot[69][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 6)
$entry = $snapshot[69][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 7)
$entry = $snapshot[69][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 8)
$entry = $snapshot[69][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 9)
$entry = $snapshot[69][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 10)
$entry = $snapshot[69][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 11)
$entry = $snapshot[69][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 12)
$entry = $snapshot[69][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 13)
$entry = $snapshot[69][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 14)
$entry = $snapshot[69][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 15)
$entry = $snapshot[69][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 16)
$entry = $snapshot[69][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(41, 17)
$entry = $snapshot[69][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 42 <- source row 44
$cell = $ws.Cells.Item(42, 1)
$entry = $snapshot[44][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 2)
$entry = $snapshot[44][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 3)
$entry = $snapshot[44][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 4)
$entry = $snapshot[44][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 5)
$entry = $snapshot[44][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 7)
$entry = $snapshot[44][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 8)
$entry = $snapshot[44][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 9)
$entry = $snapshot[44][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 10)
$entry = $snapshot[44][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 11)
$entry = $snapshot[44][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 12)
$entry = $snapshot[44][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 13)
$entry = $snapshot[44][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 14)
$entry = $snapshot[44][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 15)
$entry = $snapshot[44][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 16)
$entry = $snapshot[44][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(42, 17)
$entry = $snapshot[44][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 43 <- source row 48
$cell = $ws.Cells.Item(43, 1)
$entry = $snapshot[48][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 2)
$entry = $snapshot[48][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 3)
$entry = $snapshot[48][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 4)
$entry = $snapshot[48][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 5)
$entry = $snapshot[48][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 7)
$entry = $snapshot[48][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 8)
$entry = $snapshot[48][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 9)
$entry = $snapshot[48][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 10)
$entry = $snapshot[48][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 11)
$entry = $snapshot[48][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 12)
$entry = $snapshot[48][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 13)
$entry = $snapshot[48][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 14)
$entry = $snapshot[48][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 15)
$entry = $snapshot[48][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 16)
$entry = $snapshot[48][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(43, 17)
$entry = $snapshot[48][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 44 <- source row 67
$cell = $ws.Cells.Item(44, 1)
$entry = $snapshot[67][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 2)
$entry = $snapshot[67][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 3)
$entry = $snapshot[67][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 4)
$entry = $snapshot[67][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 5)
$entry = $snapshot[67][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 7)
$entry = $snapshot[67][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 8)
$entry = $snapshot[67][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 9)
$entry = $snapshot[67][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 10)
$entry = $snapshot[67][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 11)
$entry = $snapshot[67][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 12)
$entry = $snapshot[67][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 13)
$entry = $snapshot[67][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 14)
$entry = $snapshot[67][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 15)
$entry = $snapshot[67][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 16)
$entry = $snapshot[67][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(44, 17)
$entry = $snapshot[67][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 45 <- source row 42
$cell = $ws.Cells.Item(45, 1)
$entry = $snapshot[42][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 2)
$entry = $snapshot[42][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 3)
$entry = $snapshot[42][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 4)
$entry = $snapshot[42][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 5)
$entry = $snapshot[42][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 7)
$entry = $snapshot[42][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 8)
$entry = $snapshot[42][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 9)
$entry = $snapshot[42][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 10)
$entry = $snapshot[42][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 11)
$entry = $snapshot[42][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 12)
$entry = $snapshot[42][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 13)
$entry = $snapshot[42][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 14)
$entry = $snapshot[42][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 15)
$entry = $snapshot[42][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 16)
$entry = $snapshot[42][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(45, 17)
$entry = $snapshot[42][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 46 <- source row 38
$cell = $ws.Cells.Item(46, 1)
$entry = $snapshot[38][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 2)
$entry = $snapshot[38][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 3)
$entry = $snapshot[38][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 4)
$entry = $snapshot[38][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 5)
$entry = $snapshot[38][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 7)
$entry = $snapshot[38][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 8)
$entry = $snapshot[38][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 9)
$entry = $snapshot[38][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 10)
$entry = $snapshot[38][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 11)
$entry = $snapshot[38][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 12)
$entry = $snapshot[38][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 13)
$entry = $snapshot[38][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 14)
$entry = $snapshot[38][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 15)
$entry = $snapshot[38][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 16)
$entry = $snapshot[38][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(46, 17)
$entry = $snapshot[38][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 48 <- source row 59
$cell = $ws.Cells.Item(48, 1)
$entry = $snapshot[59][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 2)
$entry = $snapshot[59][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 3)
$entry = $snapshot[59][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 4)
$entry = $snapshot[59][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 5)
$entry = $snapshot[59][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 6)
$entry = $snapshot[59][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 7)
$entry = $snapshot[59][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 8)
$entry = $snapshot[59][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 9)
$entry = $snapshot[59][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 10)
$entry = $snapshot[59][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 11)
$entry = $snapshot[59][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 12)
$entry = $snapshot[59][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 13)
$entry = $snapshot[59][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 14)
$entry = $snapshot[59][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 15)
$entry = $snapshot[59][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 16)
$entry = $snapshot[59][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(48, 17)
$entry = $snapshot[59][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 50 <- source row 62
$cell = $ws.Cells.Item(50, 1)
$entry = $snapshot[62][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 2)
$entry = $snapshot[62][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 3)
$entry = $snapshot[62][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 4)
$entry = $snapshot[62][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 5)
$entry = $snapshot[62][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 7)
$entry = $snapshot[62][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 8)
$entry = $snapshot[62][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 9)
$entry = $snapshot[62][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 10)
$entry = $snapshot[62][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 11)
$entry = $snapshot[62][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 12)
$entry = $snapshot[62][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 13)
$entry = $snapshot[62][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 14)
$entry = $snapshot[62][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 15)
$entry = $snapshot[62][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 16)
$entry = $snapshot[62][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(50, 17)
$entry = $snapshot[62][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 51 <- source row 71
$cell = $ws.Cells.Item(51, 1)
$entry = $snapshot[71][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 2)
$entry = $snapshot[71][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 3)
$entry = $snapshot[71][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 4)
$entry = $snapshot[71][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 5)
$entry = $snapshot[71][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 7)
$entry = $snapshot[71][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 8)
$entry = $snapshot[71][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 9)
$entry = $snapshot[71][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 10)
$entry = $snapshot[71][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 11)
$entry = $snapshot[71][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 12)
$entry = $snapshot[71][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 13)
$entry = $snapshot[71][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 14)
$entry = $snapshot[71][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 15)
$entry = $snapshot[71][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 16)
$entry = $snapshot[71][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(51, 17)
$entry = $snapshot[71][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 52 <- source row 63
$cell = $ws.Cells.Item(52, 1)
$entry = $snapshot[63][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 2)
$entry = $snapshot[63][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 3)
$entry = $snapshot[63][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 4)
$entry = $snapshot[63][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 5)
$entry = $snapshot[63][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 7)
$entry = $snapshot[63][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 8)
$entry = $snapshot[63][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 9)
$entry = $snapshot[63][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 10)
$entry = $snapshot[63][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 11)
$entry = $snapshot[63][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 12)
$entry = $snapshot[63][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 13)
$entry = $snapshot[63][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 14)
$entry = $snapshot[63][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 15)
$entry = $snapshot[63][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 16)
$entry = $snapshot[63][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(52, 17)
$entry = $snapshot[63][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 53 <- source row 51
$cell = $ws.Cells.Item(53, 1)
$entry = $snapshot[51][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 2)
$entry = $snapshot[51][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 3)
$entry = $snapshot[51][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 4)
$entry = $snapshot[51][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 5)
$entry = $snapshot[51][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 7)
$entry = $snapshot[51][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 8)
$entry = $snapshot[51][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 9)
$entry = $snapshot[51][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 10)
$entry = $snapshot[51][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 11)
$entry = $snapshot[51][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 12)
$entry = $snapshot[51][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 13)
$entry = $snapshot[51][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 14)
$entry = $snapshot[51][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 15)
$entry = $snapshot[51][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 16)
$entry = $snapshot[51][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(53, 17)
$entry = $snapshot[51][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 55 <- source row 43
$cell = $ws.Cells.Item(55, 1)
$entry = $snapshot[43][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 2)
$entry = $snapshot[43][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 3)
$entry = $snapshot[43][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 4)
$entry = $snapshot[43][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 5)
$entry = $snapshot[43][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 7)
$entry = $snapshot[43][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 8)
$entry = $snapshot[43][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 9)
$entry = $snapshot[43][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 10)
$entry = $snapshot[43][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 11)
$entry = $snapshot[43][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 12)
$entry = $snapshot[43][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 13)
$entry = $snapshot[43][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 14)
$entry = $snapshot[43][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 15)
$entry = $snapshot[43][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 16)
$entry = $snapshot[43][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(55, 17)
$entry = $snapshot[43][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 56 <- source row 65
$cell = $ws.Cells.Item(56, 1)
$entry = $snapshot[65][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 2)
$entry = $snapshot[65][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 3)
$entry = $snapshot[65][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 4)
$entry = $snapshot[65][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 5)
$entry = $snapshot[65][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 7)
$entry = $snapshot[65][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 8)
$entry = $snapshot[65][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 9)
$entry = $snapshot[65][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 10)
$entry = $snapshot[65][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 11)
$entry = $snapshot[65][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 12)
$entry = $snapshot[65][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 13)
$entry = $snapshot[65][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 14)
$entry = $snapshot[65][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 15)
$entry = $snapshot[65][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 16)
$entry = $snapshot[65][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(56, 17)
$entry = $snapshot[65][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 57 <- source row 45
$cell = $ws.Cells.Item(57, 1)
$entry = $snapshot[45][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 2)
$entry = $snapshot[45][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 3)
$entry = $snapshot[45][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 4)
$entry = $snapshot[45][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 5)
$entry = $snapshot[45][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 7)
$entry = $snapshot[45][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 8)
$entry = $snapshot[45][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 9)
$entry = $snapshot[45][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 10)
$entry = $snapshot[45][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 11)
$entry = $snapshot[45][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 12)
$entry = $snapshot[45][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 13)
$entry = $snapshot[45][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 14)
$entry = $snapshot[45][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 15)
$entry = $snapshot[45][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 16)
$entry = $snapshot[45][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(57, 17)
$entry = $snapshot[45][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(57, 6).ClearContents()

# Target row 58 <- source row 35
$cell = $ws.Cells.Item(58, 1)
$entry = $snapshot[35][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 2)
$entry = $snapshot[35][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 3)
$entry = $snapshot[35][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 4)
$entry = $snapshot[35][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 5)
$entry = $snapshot[35][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 6)
$entry = $snapshot[35][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 7)
$entry = $snapshot[35][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 8)
$entry = $snapshot[35][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 9)
$entry = $snapshot[35][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 10)
$entry = $snapshot[35][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 11)
$entry = $snapshot[35][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 12)
$entry = $snapshot[35][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 13)
$entry = $snapshot[35][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 14)
$entry = $snapshot[35][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 15)
$entry = $snapshot[35][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 16)
$entry = $snapshot[35][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(58, 17)
$entry = $snapshot[35][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 59 <- source row 40
$cell = $ws.Cells.Item(59, 1)
$entry = $snapshot[40][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 2)
$entry = $snapshot[40][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 3)
$entry = $snapshot[40][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 4)
$entry = $snapshot[40][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 5)
$entry = $snapshot[40][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 7)
$entry = $snapshot[40][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 8)
$entry = $snapshot[40][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 9)
$entry = $snapshot[40][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 10)
$entry = $snapshot[40][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 11)
$entry = $snapshot[40][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 12)
$entry = $snapshot[40][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 13)
$entry = $snapshot[40][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 14)
$entry = $snapshot[40][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 15)
$entry = $snapshot[40][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 16)
$entry = $snapshot[40][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(59, 17)
$entry = $snapshot[40][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(59, 6).ClearContents()

# Target row 60 <- source row 34
$cell = $ws.Cells.Item(60, 1)
$entry = $snapshot[34][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 2)
$entry = $snapshot[34][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 3)
$entry = $snapshot[34][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 4)
$entry = $snapshot[34][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 5)
$entry = $snapshot[34][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 7)
$entry = $snapshot[34][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 8)
$entry = $snapshot[34][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 9)
$entry = $snapshot[34][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 10)
$entry = $snapshot[34][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 11)
$entry = $snapshot[34][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 12)
$entry = $snapshot[34][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 13)
$entry = $snapshot[34][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 14)
$entry = $snapshot[34][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 15)
$entry = $snapshot[34][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 16)
$entry = $snapshot[34][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(60, 17)
$entry = $snapshot[34][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 61 <- source row 36
$cell = $ws.Cells.Item(61, 1)
$entry = $snapshot[36][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 2)
$entry = $snapshot[36][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 3)
$entry = $snapshot[36][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 4)
$entry = $snapshot[36][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 5)
$entry = $snapshot[36][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 7)
$entry = $snapshot[36][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 8)
$entry = $snapshot[36][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 9)
$entry = $snapshot[36][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 10)
$entry = $snapshot[36][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 11)
$entry = $snapshot[36][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 12)
$entry = $snapshot[36][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 13)
$entry = $snapshot[36][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 14)
$entry = $snapshot[36][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 15)
$entry = $snapshot[36][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 16)
$entry = $snapshot[36][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(61, 17)
$entry = $snapshot[36][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 62 <- source row 41
$cell = $ws.Cells.Item(62, 1)
$entry = $snapshot[41][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 2)
$entry = $snapshot[41][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 3)
$entry = $snapshot[41][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 4)
$entry = $snapshot[41][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 5)
$entry = $snapshot[41][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 7)
$entry = $snapshot[41][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 8)
$entry = $snapshot[41][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 9)
$entry = $snapshot[41][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 10)
$entry = $snapshot[41][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 11)
$entry = $snapshot[41][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 12)
$entry = $snapshot[41][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 13)
$entry = $snapshot[41][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 14)
$entry = $snapshot[41][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 15)
$entry = $snapshot[41][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 16)
$entry = $snapshot[41][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(62, 17)
$entry = $snapshot[41][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 63 <- source row 61
$cell = $ws.Cells.Item(63, 1)
$entry = $snapshot[61][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 2)
$entry = $snapshot[61][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 3)
$entry = $snapshot[61][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 4)
$entry = $snapshot[61][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 5)
$entry = $snapshot[61][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 7)
$entry = $snapshot[61][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 8)
$entry = $snapshot[61][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 9)
$entry = $snapshot[61][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 10)
$entry = $snapshot[61][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 11)
$entry = $snapshot[61][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 12)
$entry = $snapshot[61][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 13)
$entry = $snapshot[61][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 14)
$entry = $snapshot[61][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 15)
$entry = $snapshot[61][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 16)
$entry = $snapshot[61][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(63, 17)
$entry = $snapshot[61][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 64 <- source row 66
$cell = $ws.Cells.Item(64, 1)
$entry = $snapshot[66][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 2)
$entry = $snapshot[66][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 3)
$entry = $snapshot[66][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 4)
$entry = $snapshot[66][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 5)
$entry = $snapshot[66][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 6)
$entry = $snapshot[66][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 7)
$entry = $snapshot[66][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 8)
$entry = $snapshot[66][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 9)
$entry = $snapshot[66][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 10)
$entry = $snapshot[66][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 11)
$entry = $snapshot[66][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 12)
$entry = $snapshot[66][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 13)
$entry = $snapshot[66][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 14)
$entry = $snapshot[66][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 15)
$entry = $snapshot[66][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 16)
$entry = $snapshot[66][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(64, 17)
$entry = $snapshot[66][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 65 <- source row 64
$cell = $ws.Cells.Item(65, 1)
$entry = $snapshot[64][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 2)
$entry = $snapshot[64][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 3)
$entry = $snapshot[64][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 4)
$entry = $snapshot[64][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 5)
$entry = $snapshot[64][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 7)
$entry = $snapshot[64][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 8)
$entry = $snapshot[64][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 9)
$entry = $snapshot[64][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 10)
$entry = $snapshot[64][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 11)
$entry = $snapshot[64][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 12)
$entry = $snapshot[64][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 13)
$entry = $snapshot[64][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 14)
$entry = $snapshot[64][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 15)
$entry = $snapshot[64][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 16)
$entry = $snapshot[64][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(65, 17)
$entry = $snapshot[64][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 66 <- source row 68
$cell = $ws.Cells.Item(66, 1)
$entry = $snapshot[68][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 2)
$entry = $snapshot[68][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 3)
$entry = $snapshot[68][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 4)
$entry = $snapshot[68][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 5)
$entry = $snapshot[68][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 7)
$entry = $snapshot[68][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 8)
$entry = $snapshot[68][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 9)
$entry = $snapshot[68][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 10)
$entry = $snapshot[68][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 11)
$entry = $snapshot[68][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 12)
$entry = $snapshot[68][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 13)
$entry = $snapshot[68][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 14)
$entry = $snapshot[68][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 15)
$entry = $snapshot[68][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 16)
$entry = $snapshot[68][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(66, 17)
$entry = $snapshot[68][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(66, 6).ClearContents()

# Target row 67 <- source row 58
$cell = $ws.Cells.Item(67, 1)
$entry = $snapshot[58][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 2)
$entry = $snapshot[58][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 3)
$entry = $snapshot[58][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 4)
$entry = $snapshot[58][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 5)
$entry = $snapshot[58][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 6)
$entry = $snapshot[58][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 7)
$entry = $snapshot[58][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 8)
$entry = $snapshot[58][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 9)
$entry = $snapshot[58][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 10)
$entry = $snapshot[58][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 11)
$entry = $snapshot[58][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 12)
$entry = $snapshot[58][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 13)
$entry = $snapshot[58][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 14)
$entry = $snapshot[58][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 15)
$entry = $snapshot[58][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 16)
$entry = $snapshot[58][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(67, 17)
$entry = $snapshot[58][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 68 <- source row 37
$cell = $ws.Cells.Item(68, 1)
$entry = $snapshot[37][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 2)
$entry = $snapshot[37][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 3)
$entry = $snapshot[37][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 4)
$entry = $snapshot[37][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 5)
$entry = $snapshot[37][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 7)
$entry = $snapshot[37][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 8)
$entry = $snapshot[37][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 9)
$entry = $snapshot[37][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 10)
$entry = $snapshot[37][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 11)
$entry = $snapshot[37][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 12)
$entry = $snapshot[37][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 13)
$entry = $snapshot[37][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 14)
$entry = $snapshot[37][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 15)
$entry = $snapshot[37][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 16)
$entry = $snapshot[37][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(68, 17)
$entry = $snapshot[37][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 69 <- source row 50
$cell = $ws.Cells.Item(69, 1)
$entry = $snapshot[50][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 2)
$entry = $snapshot[50][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 3)
$entry = $snapshot[50][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 4)
$entry = $snapshot[50][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 5)
$entry = $snapshot[50][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 7)
$entry = $snapshot[50][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 8)
$entry = $snapshot[50][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 9)
$entry = $snapshot[50][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 10)
$entry = $snapshot[50][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 11)
$entry = $snapshot[50][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 12)
$entry = $snapshot[50][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 13)
$entry = $snapshot[50][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 14)
$entry = $snapshot[50][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 15)
$entry = $snapshot[50][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 16)
$entry = $snapshot[50][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(69, 17)
$entry = $snapshot[50][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$ws.Cells.Item(69, 6).ClearContents()

# Target row 70 <- source row 57
$cell = $ws.Cells.Item(70, 1)
$entry = $snapshot[57][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 2)
$entry = $snapshot[57][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 3)
$entry = $snapshot[57][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 4)
$entry = $snapshot[57][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 5)
$entry = $snapshot[57][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 6)
$entry = $snapshot[57][6]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 7)
$entry = $snapshot[57][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 8)
$entry = $snapshot[57][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 9)
$entry = $snapshot[57][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 10)
$entry = $snapshot[57][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 11)
$entry = $snapshot[57][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 12)
$entry = $snapshot[57][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 13)
$entry = $snapshot[57][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 14)
$entry = $snapshot[57][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 15)
$entry = $snapshot[57][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 16)
$entry = $snapshot[57][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(70, 17)
$entry = $snapshot[57][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 71 <- source row 60
$cell = $ws.Cells.Item(71, 1)
$entry = $snapshot[60][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 2)
$entry = $snapshot[60][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 3)
$entry = $snapshot[60][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 4)
$entry = $snapshot[60][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 5)
$entry = $snapshot[60][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 7)
$entry = $snapshot[60][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 8)
$entry = $snapshot[60][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 9)
$entry = $snapshot[60][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 10)
$entry = $snapshot[60][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 11)
$entry = $snapshot[60][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 12)
$entry = $snapshot[60][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 13)
$entry = $snapshot[60][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 14)
$entry = $snapshot[60][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 15)
$entry = $snapshot[60][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 16)
$entry = $snapshot[60][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(71, 17)
$entry = $snapshot[60][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 72 <- source row 39
$cell = $ws.Cells.Item(72, 1)
$entry = $snapshot[39][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 2)
$entry = $snapshot[39][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 3)
$entry = $snapshot[39][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 4)
$entry = $snapshot[39][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 5)
$entry = $snapshot[39][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 7)
$entry = $snapshot[39][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 8)
$entry = $snapshot[39][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 9)
$entry = $snapshot[39][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 10)
$entry = $snapshot[39][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 11)
$entry = $snapshot[39][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 12)
$entry = $snapshot[39][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 13)
$entry = $snapshot[39][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 14)
$entry = $snapshot[39][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 15)
$entry = $snapshot[39][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 16)
$entry = $snapshot[39][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(72, 17)
$entry = $snapshot[39][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 73 <- source row 52
$cell = $ws.Cells.Item(73, 1)
$entry = $snapshot[52][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 2)
$entry = $snapshot[52][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 3)
$entry = $snapshot[52][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 4)
$entry = $snapshot[52][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 5)
$entry = $snapshot[52][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 7)
$entry = $snapshot[52][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 8)
$entry = $snapshot[52][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 9)
$entry = $snapshot[52][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 10)
$entry = $snapshot[52][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 11)
$entry = $snapshot[52][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 12)
$entry = $snapshot[52][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 13)
$entry = $snapshot[52][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 14)
$entry = $snapshot[52][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 15)
$entry = $snapshot[52][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 16)
$entry = $snapshot[52][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(73, 17)
$entry = $snapshot[52][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 74 <- source row 75
$cell = $ws.Cells.Item(74, 1)
$entry = $snapshot[75][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 2)
$entry = $snapshot[75][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 3)
$entry = $snapshot[75][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 4)
$entry = $snapshot[75][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 5)
$entry = $snapshot[75][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 7)
$entry = $snapshot[75][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 8)
$entry = $snapshot[75][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 9)
$entry = $snapshot[75][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 10)
$entry = $snapshot[75][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 11)
$entry = $snapshot[75][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 12)
$entry = $snapshot[75][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 13)
$entry = $snapshot[75][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 14)
$entry = $snapshot[75][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 15)
$entry = $snapshot[75][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 16)
$entry = $snapshot[75][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(74, 17)
$entry = $snapshot[75][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Target row 75 <- source row 74
$cell = $ws.Cells.Item(75, 1)
$entry = $snapshot[74][1]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 2)
$entry = $snapshot[74][2]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 3)
$entry = $snapshot[74][3]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 4)
$entry = $snapshot[74][4]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 5)
$entry = $snapshot[74][5]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 7)
$entry = $snapshot[74][7]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 8)
$entry = $snapshot[74][8]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 9)
$entry = $snapshot[74][9]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 10)
$entry = $snapshot[74][10]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 11)
$entry = $snapshot[74][11]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 12)
$entry = $snapshot[74][12]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 13)
$entry = $snapshot[74][13]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 14)
$entry = $snapshot[74][14]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 15)
$entry = $snapshot[74][15]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 16)
$entry = $snapshot[74][16]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }
$cell = $ws.Cells.Item(75, 17)
$entry = $snapshot[74][17]
if ($entry.Kind -eq "Formula") { $cell.Formula = $entry.Val } else { $cell.Value2 = $entry.Val }

# Bump column C (Förändrad) by 1 day for every data row 2..75
for ($r = 2; $r -le 75; $r++) {
    $c = $ws.Cells.Item($r, 3)
    $c.Value2 = $c.Value2 + 1
}